$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '60.026.77'
$ws.Range('E2').Value = '  -6.32%  '
$ws.Range('D3').Value = '3.278.14'
$ws.Range('E3').Value = '  -5.61%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.46%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.279.57'
$ws.Range('E8').Value = '  -5.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.370'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.52%  '
$ws.Range('D13').Value = '3.836.97'
$ws.Range('E13').Value = '  -5.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.119'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.280.29'
$ws.Range('E15').Value = '  -5.52%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000167'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '60.249.99'
$ws.Range('E18').Value = '  -5.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('E21').Value = '  -9.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '350.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.550'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '3.405.18'
$ws.Range('E25').Value = '  -5.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.54%  '
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.20%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').Value = '3.306.15'
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.75'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '157.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0748'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.62%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.736'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.60%  '
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.857'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.00%  '
